$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 266.85715
$ws.Range("I53").Value = 338.9
$ws.Range("K53").Value = 338.9
$ws.Range("M53").Value = 298.1

$ws.Range("H76").Value = 1981
$ws.Range("I76").Value = 1981
$ws.Range("K76").Value = 1981
$ws.Range("M76").Value = -1666

$ws.Range("H79").Value = 1981
$ws.Range("I79").Value = 1981
$ws.Range("K79").Value = 1981
$ws.Range("M79").Value = -889

$ws.Range("H97").Value = 2036.7778
$ws.Range("J97").Value = 2036.7778
$ws.Range("L97").Value = 6110.3334
$ws.Range("N97").Value = -7102.3334

$ws.Range("H99").Value = 2347362.8
$ws.Range("I99").Value = 3281308
$ws.Range("K99").Value = 9843924
$ws.Range("M99").Value = -9842426

$ws.Range("H101").Value = 7519975.5
$ws.Range("J101").Value = 925.875
$ws.Range("L101").Value = 2777.625
$ws.Range("N101").Value = -6021.625

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("N121").Value = 0

$ws.Range("H126").Value = 84778
$ws.Range("J126").Value = 84778
$ws.Range("L126").Value = 84778
$ws.Range("N126").Value = -94658

$ws.Range("H132").Value = 4085.6
$ws.Range("I132").Value = 4396.9287
$ws.Range("K132").Value = 13190.7861
$ws.Range("M132").Value = -10660.7861

$ws.Range("H135").Value = 5700.0454
$ws.Range("I135").Value = 5915.8423
$ws.Range("J135").Value = 4333.3335
$ws.Range("K135").Value = 53242.58070000001
$ws.Range("L135").Value = 39000.0015
$ws.Range("M135").Value = -50707.58070000001
$ws.Range("N135").Value = -44070.0015

$ws.Range("H137").Value = 668963.6
$ws.Range("I137").Value = 1176228.2
$ws.Range("J137").Value = 16766.285
$ws.Range("K137").Value = 3528684.6
$ws.Range("L137").Value = 50298.855
$ws.Range("M137").Value = -3526134.6
$ws.Range("N137").Value = -55398.855


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 304.2
$ws.Range("I5").Value = 304.2
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 304.2
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -192.2

$ws.Range("H32").Value = 3180.054
$ws.Range("I32").Value = 3133.2
$ws.Range("J32").Value = 4000
$ws.Range("K32").Value = 3133.2
$ws.Range("L32").Value = 4000
$ws.Range("M32").Value = -2846.2
$ws.Range("N32").Value = -4574

$ws.Range("H45").Value = 169635.47
$ws.Range("I45").Value = 270438.88
$ws.Range("K45").Value = 270438.88
$ws.Range("M45").Value = -270061.88

$ws.Range("H132").Value = 3237.6453
$ws.Range("I132").Value = 2431.1
$ws.Range("K132").Value = 7293.299999999999
$ws.Range("M132").Value = -4763.299999999999

$ws.Range("H139").Value = 76357.39999999999
$ws.Range("J139").Value = 76357.39999999999
$ws.Range("L139").Value = 76357.39999999999
$ws.Range("N139").Value = -86637.39999999999


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 304.2
$ws.Range("I4").Value = 304.2
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 304.2
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -189.2


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2836.2666
$ws.Range("I16").Value = 2977.2727
$ws.Range("J16").Value = 2448.5
$ws.Range("K16").Value = 2977.2727
$ws.Range("L16").Value = 2448.5
$ws.Range("M16").Value = -2690.2727
$ws.Range("N16").Value = -3022.5

$ws.Range("H22").Value = 693.3570999999999
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 650.7
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 650.7
$ws.Range("M22").Value = -450
$ws.Range("N22").Value = -1350.7

$ws.Range("H113").Value = 2836.2666
$ws.Range("I113").Value = 2977.2727
$ws.Range("J113").Value = 2448.5
$ws.Range("K113").Value = 2977.2727
$ws.Range("L113").Value = 2448.5
$ws.Range("M113").Value = -807.2727
$ws.Range("N113").Value = -6788.5

$ws.Range("H132").Value = 18547.594
$ws.Range("I132").Value = 6193.7144
$ws.Range("K132").Value = 18581.1432
$ws.Range("M132").Value = -16051.1432

$ws.Range("H134").Value = 4820032
$ws.Range("I134").Value = 7829065
$ws.Range("K134").Value = 23487195
$ws.Range("M134").Value = -23484660


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 41.2
$ws.Range("J12").Value = 35.333332
$ws.Range("L12").Value = 105.999996
$ws.Range("N12").Value = -451.999996

$ws.Range("H23").Value = 178.41667
$ws.Range("J23").Value = 238.57143
$ws.Range("L23").Value = 715.71429
$ws.Range("N23").Value = -1185.71429

$ws.Range("H46").Value = 2208.0908
$ws.Range("J46").Value = 5199.75
$ws.Range("L46").Value = 15599.25
$ws.Range("N46").Value = -15781.25

$ws.Range("H108").Value = 1722.5834
$ws.Range("I108").Value = 908
$ws.Range("K108").Value = 2724
$ws.Range("M108").Value = 156

$ws.Range("H114").Value = 913.5
$ws.Range("I114").Value = 875
$ws.Range("J114").Value = 932.75
$ws.Range("K114").Value = 2625
$ws.Range("L114").Value = 2798.25
$ws.Range("M114").Value = 629
$ws.Range("N114").Value = -9306.25

$ws.Range("H132").Value = 56809.445
$ws.Range("I132").Value = 900
$ws.Range("J132").Value = 63798.125
$ws.Range("K132").Value = 8100
$ws.Range("L132").Value = 574183.125
$ws.Range("M132").Value = -5570
$ws.Range("N132").Value = -579243.125

$ws.Range("H137").Value = 3121.4614
$ws.Range("J137").Value = 8331.666999999999
$ws.Range("L137").Value = 24995.001
$ws.Range("N137").Value = -35195.001


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 9023.117
$ws.Range("I102").Value = 10242.357
$ws.Range("J102").Value = 3333.3333
$ws.Range("K102").Value = 10242.357
$ws.Range("L102").Value = 3333.3333
$ws.Range("M102").Value = -8620.357
$ws.Range("N102").Value = -6577.3333

$ws.Range("H113").Value = 3299.889
$ws.Range("I113").Value = 2133.3333
$ws.Range("J113").Value = 3883.1667
$ws.Range("K113").Value = 2133.3333
$ws.Range("L113").Value = 3883.1667
$ws.Range("M113").Value = 36.66670000000022
$ws.Range("N113").Value = -8223.1667

$ws.Range("H126").Value = 19033
$ws.Range("I126").Value = 32801.2
$ws.Range("J126").Value = 12774.728
$ws.Range("K126").Value = 98403.59999999999
$ws.Range("L126").Value = 38324.18399999999
$ws.Range("M126").Value = -95933.59999999999
$ws.Range("N126").Value = -43264.18399999999

$ws.Range("H132").Value = 3432.45
$ws.Range("I132").Value = 2058.25
$ws.Range("K132").Value = 6174.75
$ws.Range("M132").Value = -3644.75

$ws.Range("H136").Value = 124050.664
$ws.Range("J136").Value = 124050.664
$ws.Range("L136").Value = 372151.992
$ws.Range("N136").Value = -377251.992

$ws.Range("H139").Value = 89999.5
$ws.Range("J139").Value = 89999.5
$ws.Range("L139").Value = 89999.5
$ws.Range("N139").Value = -100279.5


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 27580.857
$ws.Range("I7").Value = 28484.9
$ws.Range("J7").Value = 9500
$ws.Range("K7").Value = 28484.9
$ws.Range("L7").Value = 9500
$ws.Range("M7").Value = -28372.9
$ws.Range("N7").Value = -9724

$ws.Range("H40").Value = 32249.938
$ws.Range("I40").Value = 40273.184
$ws.Range("K40").Value = 40273.184
$ws.Range("M40").Value = -40137.184

$ws.Range("H61").Value = 3876
$ws.Range("I61").Value = 2752
$ws.Range("K61").Value = 2752
$ws.Range("M61").Value = -2550

$ws.Range("H93").Value = 3356.8333
$ws.Range("I93").Value = 4282.875
$ws.Range("K93").Value = 4282.875
$ws.Range("M93").Value = -3034.875

$ws.Range("H113").Value = 3876
$ws.Range("I113").Value = 2752
$ws.Range("K113").Value = 2752
$ws.Range("M113").Value = -582

$ws.Range("H122").Value = 6250
$ws.Range("J122").Value = 8000
$ws.Range("L122").Value = 24000
$ws.Range("N122").Value = -28900

$ws.Range("H123").Value = 68999.5
$ws.Range("J123").Value = 68999.5
$ws.Range("L123").Value = 68999.5
$ws.Range("N123").Value = -78799.5

$ws.Range("H126").Value = 27580.857
$ws.Range("I126").Value = 28484.9
$ws.Range("J126").Value = 9500
$ws.Range("K126").Value = 85454.70000000001
$ws.Range("L126").Value = 28500
$ws.Range("M126").Value = -82984.70000000001
$ws.Range("N126").Value = -33440

$ws.Range("H132").Value = 4585437.5
$ws.Range("I132").Value = 6876131
$ws.Range("J132").Value = 4050
$ws.Range("K132").Value = 20628393
$ws.Range("L132").Value = 12150
$ws.Range("M132").Value = -20625863
$ws.Range("N132").Value = -17210

$ws.Range("H136").Value = 13342.723
$ws.Range("J136").Value = 9874.166999999999
$ws.Range("L136").Value = 29622.501
$ws.Range("N136").Value = -34722.501


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 40000
$ws.Range("J119").Value = 40000
$ws.Range("L119").Value = 40000
$ws.Range("N119").Value = -49676

$ws.Range("H122").Value = 37516.707
$ws.Range("I122").Value = 4099.875
$ws.Range("K122").Value = 12299.625
$ws.Range("M122").Value = -9849.625

$ws.Range("H126").Value = 23832.35
$ws.Range("I126").Value = 28950.867
$ws.Range("K126").Value = 86852.601
$ws.Range("M126").Value = -84382.601

$ws.Range("H132").Value = 22454.842
$ws.Range("I132").Value = 23332.47
$ws.Range("K132").Value = 69997.41
$ws.Range("M132").Value = -67467.41

